$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fecundity rows (152-266) captured in this batch of processing
$data = New-Object "object[,]" 115,5

$data[0,0] = "d§JBCHY1/1/50§Q§273"
$data[0,1] = 15.05
$data[0,2] = "n"
$data[0,4] = "y"

$data[1,0] = "b§KVEDG1/1/80§Q§279"
$data[1,1] = 15.37
$data[1,2] = "n"
$data[1,4] = "y"

$data[2,0] = "e§VRCAN/1/40§Q§504"
$data[2,1] = 15.3
$data[2,2] = "n"
$data[2,4] = "y"

$data[3,0] = "b§WSSWM3/1/0§Q§569"
$data[3,1] = 14.33
$data[3,2] = "n"
$data[3,4] = "y"

$data[4,0] = "b§JBCHY1/1/50§Q§227"
$data[4,1] = 14.48
$data[4,2] = "n"
$data[4,4] = "y"

$data[5,0] = "a§DVGM/20§Q§105"
$data[5,1] = 16.13
$data[5,2] = "n"
$data[5,4] = "y"

$data[6,0] = "a§KVEDG1/1/80§Q§275"
$data[6,1] = 16.19
$data[6,2] = "n"
$data[6,4] = "y"

$data[7,0] = "a§DVGM/20§Q§107"
$data[7,1] = 15.85
$data[7,2] = "n"
$data[7,4] = "y"

$data[8,0] = "d§CRSOSO/3/40§Q§62"
$data[8,1] = 15.12
$data[8,2] = "n"
$data[8,4] = "y"

$data[9,0] = "d§CBMCK1/3/60§N§i115"
$data[9,1] = 15.43
$data[9,2] = "n"
$data[9,4] = "y"

$data[10,0] = "d§KVEDG1/1/80§Q§294"
$data[10,1] = 15.44
$data[10,2] = "n"
$data[10,4] = "y"

$data[11,0] = "c§CBMCK1/1/0§N§42"
$data[11,1] = 16.35
$data[11,2] = "n"
$data[11,4] = "y"

$data[12,0] = "d§VRCAN/1/40§Q§501"
$data[12,1] = 17.3
$data[12,2] = "n"
$data[12,4] = "y"

$data[13,0] = "d§JARI1/1/80§Q§180"
$data[13,1] = 15.1
$data[13,2] = "n"
$data[13,4] = "y"

$data[14,0] = "b§JBCHY1/1/50§Q§242"
$data[14,1] = 20.42
$data[14,2] = "n"
$data[14,3] = "y"
$data[14,4] = "y"

$data[15,0] = "b§JBCHY1/1/50§Q§242"
$data[15,1] = 6.22
$data[15,2] = "y"
$data[15,3] = "y"
$data[15,4] = "n"

$data[16,0] = "b§MHBUR1/4/20§Q§343"
$data[16,1] = 14.46
$data[16,2] = "n"
$data[16,4] = "y"

$data[17,0] = "c§CWRIC2/X2§Q§91"
$data[17,1] = 14.64
$data[17,2] = "n"
$data[17,4] = "y"

$data[18,0] = "c§JBBLB2/1/20§Q§210"
$data[18,1] = 15.04
$data[18,2] = "n"
$data[18,4] = "y"

$data[19,0] = "c§VRCAN/1/40§Q§491"
$data[19,1] = 14.45
$data[19,2] = "n"
$data[19,4] = "y"

$data[20,0] = "d§JARI1/1/80§Q§184"
$data[20,1] = 6.33
$data[20,2] = "y"
$data[20,4] = "n"

$data[21,0] = "a§SMITH1/20§Q§456"
$data[21,1] = 6.6
$data[21,2] = "y"
$data[21,4] = "n"

$data[22,0] = "b§VRPET2/3/20§Q§510"
$data[22,1] = 6.57
$data[22,2] = "y"
$data[22,4] = "n"

$data[23,0] = "maple/cnt§11"
$data[23,1] = 7.37
$data[23,2] = "y"
$data[23,4] = "n"

$data[24,0] = "a§VSGARO1/1/0§Q§533"
$data[24,1] = 7.17
$data[24,2] = "y"
$data[24,4] = "n"

$data[25,0] = "d§MHBUR1/4/20§Q§350"
$data[25,1] = 5.72
$data[25,2] = "y"
$data[25,4] = "n"

$data[26,0] = "e§RULEB1/20§N§429"
$data[26,1] = 6.12
$data[26,2] = "y"
$data[26,4] = "n"

$data[27,0] = "e§WSSWM3/1/0§Q§575"
$data[27,1] = 6.33
$data[27,2] = "y"
$data[27,4] = "n"

$data[28,0] = "maple/cnt§27"
$data[28,1] = 6.75
$data[28,2] = "y"
$data[28,4] = "n"

$data[29,0] = "maple/cnt§15"
$data[29,1] = 6.9
$data[29,2] = "y"
$data[29,4] = "n"

$data[30,0] = "a§CBMCK1/1/0§N§36"
$data[30,1] = 6.13
$data[30,2] = "y"
$data[30,4] = "n"

$data[31,0] = "c§CWRIC2/X2§Q§92"
$data[31,1] = 6.79
$data[31,2] = "y"
$data[31,4] = "n"

$data[32,0] = "b§VRCAN/1/40§Q§482"
$data[32,1] = 6.59
$data[32,2] = "y"
$data[32,4] = "n"

$data[33,0] = "d§CWRIC2/X2§Q§98"
$data[33,1] = 6.14
$data[33,2] = "y"
$data[33,4] = "n"

$data[34,0] = "b§SMITH1/20§Q§464"
$data[34,1] = 6.04
$data[34,2] = "y"
$data[34,4] = "n"

$data[35,0] = "d§JWBOY/2/80§Q§262"
$data[35,1] = 6.21
$data[35,2] = "y"
$data[35,4] = "n"

$data[36,0] = "e§EFCC2/4/80§Q§157"
$data[36,1] = 6.06
$data[36,2] = "y"
$data[36,4] = "n"

$data[37,0] = "d§CWRIC2/X2§Q§96"
$data[37,1] = 6.21
$data[37,2] = "y"
$data[37,4] = "n"

$data[38,0] = "a§CBMCK1/1/0§N§35"
$data[38,1] = 6.38
$data[38,2] = "y"
$data[38,4] = "n"

$data[39,0] = "d§DVGM/20§Q§122"
$data[39,1] = 6.1
$data[39,2] = "y"
$data[39,4] = "n"

$data[40,0] = "d§MSMID1/1/20§Q§392"
$data[40,1] = 6.45
$data[40,2] = "y"
$data[40,4] = "n"

$data[41,0] = "a§WSSWM3/1/0§Q§555"
$data[41,1] = 6.56
$data[41,2] = "y"
$data[41,4] = "n"

$data[42,0] = "a§SMAKC1/B§N§432"
$data[42,1] = 6.33
$data[42,2] = "y"
$data[42,4] = "n"

$data[43,0] = "d§MAVBEL2/1/20§N§325"
$data[43,1] = 6.48
$data[43,2] = "y"
$data[43,4] = "n"

$data[44,0] = "b§EFCC2/4/80§Q§144"
$data[44,1] = 6.28
$data[44,2] = "y"
$data[44,4] = "n"

$data[45,0] = "c§MHBUR1/4/20§Q§344"
$data[45,1] = 6.11
$data[45,2] = "y"
$data[45,4] = "n"

$data[46,0] = "d§VRPET2/3/20§Q§506"
$data[46,1] = 5.95
$data[46,2] = "y"
$data[46,4] = "n"

$data[47,0] = "b§JARI1/1/80§Q§172"
$data[47,1] = 6.52
$data[47,2] = "y"
$data[47,4] = "n"

$data[48,0] = "e§WSSWM3/1/0§Q§573"
$data[48,1] = 6.38
$data[48,2] = "y"
$data[48,4] = "n"

$data[49,0] = "e§MAVBEL2/1/20§N§329"
$data[49,1] = 5.96
$data[49,2] = "y"
$data[49,4] = "n"

$data[50,0] = "b§JBBLB2/1/20§Q§202"
$data[50,1] = 6.89
$data[50,2] = "y"
$data[50,4] = "n"

$data[51,0] = "maple/cnt§14"
$data[51,1] = 6.99
$data[51,2] = "y"
$data[51,4] = "n"

$data[52,0] = "e§PDVRT1/20§Q§418"
$data[52,1] = 6.39
$data[52,2] = "y"
$data[52,4] = "n"

$data[53,0] = "b§MHNAT1/2/0§Q§358"
$data[53,1] = 6.36
$data[53,2] = "y"
$data[53,4] = "n"

$data[54,0] = "b§VRCAN/1/40§Q§483"
$data[54,1] = 6.33
$data[54,2] = "y"
$data[54,4] = "n"

$data[55,0] = "b§WSSWM3/1/0§Q§560"
$data[55,1] = 6.12
$data[55,2] = "y"
$data[55,4] = "n"

$data[56,0] = "e§JWBOY/2/80§Q§269"
$data[56,1] = 6.4
$data[56,2] = "y"
$data[56,4] = "n"

$data[57,0] = "e§JWBOY/2/80§Q§268"
$data[57,1] = 6.33
$data[57,2] = "y"
$data[57,4] = "n"

$data[58,0] = "b§JARI1/1/80§Q§169"
$data[58,1] = 6.66
$data[58,2] = "y"
$data[58,4] = "n"

$data[59,0] = "d§CRSOSO/3/40§Q§65"
$data[59,1] = 6.44
$data[59,2] = "y"
$data[59,4] = "n"

$data[60,0] = "c§KVEDG1/1/80§Q§283"
$data[60,1] = 7.6
$data[60,2] = "y"
$data[60,4] = "n"

$data[61,0] = "c§PDVRT1/20§Q§409"
$data[61,1] = 6.42
$data[61,2] = "y"
$data[61,4] = "n"

$data[62,0] = "b|CRSOSO-3-40|Q|76"
$data[62,1] = 6.3
$data[62,2] = "y"
$data[62,4] = "n"

$data[63,0] = "maple/cnt§21"
$data[63,1] = 6.42
$data[63,2] = "y"
$data[63,4] = "n"

$data[64,0] = "a§CWRIC2/X2§Q§84"
$data[64,1] = 5.99
$data[64,2] = "y"
$data[64,4] = "n"

$data[65,0] = "a§MHBUR1/4/20§Q§338"
$data[65,1] = 6.37
$data[65,2] = "y"
$data[65,4] = "n"

$data[66,0] = "d§VSGARO1/1/0§Q§540"
$data[66,1] = 7.21
$data[66,2] = "y"
$data[66,4] = "n"

$data[67,0] = "c§DVGM/20§Q§115"
$data[67,1] = 6.33
$data[67,2] = "y"
$data[67,4] = "n"

$data[68,0] = "maple/cnt§19"
$data[68,1] = 6.35
$data[68,2] = "y"
$data[68,4] = "n"

$data[69,0] = "e§JARI1/1/80§Q§189"
$data[69,1] = 6.39
$data[69,2] = "y"
$data[69,4] = "n"

$data[70,0] = "maple/cnt§20"
$data[70,1] = 7.34
$data[70,2] = "y"
$data[70,4] = "n"

$data[71,0] = "maple/cnt§23"
$data[71,1] = 7.08
$data[71,2] = "y"
$data[71,4] = "n"

$data[72,0] = "e§MSMID1/1/20§Q§602"
$data[72,1] = 6.73
$data[72,2] = "y"
$data[72,4] = "n"

$data[73,0] = "c§MHBUR1/4/20§Q§336"
$data[73,1] = 6.67
$data[73,2] = "y"
$data[73,4] = "n"

$data[74,0] = "b§MHNAT1/2/0§Q§359"
$data[74,1] = 6.5
$data[74,2] = "y"
$data[74,4] = "n"

$data[75,0] = "d§BWPEM1/9/0§Q§18"
$data[75,1] = 9.01
$data[75,2] = "n"
$data[75,4] = "n"

$data[76,0] = "a§MSMID1/1/20§Q§366"
$data[76,1] = 9.92
$data[76,2] = "n"
$data[76,4] = "n"

$data[77,0] = "a§KVEDG1/1/80§Q§273"
$data[77,1] = 10.01
$data[77,2] = "n"
$data[77,4] = "n"

$data[78,0] = "a|CWRIC2-X2|Q|80"
$data[78,1] = 8.64
$data[78,2] = "n"
$data[78,3] = "y"
$data[78,4] = "n"

$data[79,0] = "a|CWRIC2-X2|Q|80"
$data[79,1] = 6.46
$data[79,2] = "y"
$data[79,3] = "y"
$data[79,4] = "n"

$data[80,0] = "e§DVGM/20§Q§132"
$data[80,1] = 18.96
$data[80,2] = "n"
$data[80,4] = "n"

$data[81,0] = "b§JBCHY1/1/50§Q§223"
$data[81,1] = 9.47
$data[81,2] = "n"
$data[81,4] = "n"

$data[82,0] = "a§SEBRN/1/20§N§i145"
$data[82,1] = 19.37
$data[82,2] = "n"
$data[82,4] = "n"

$data[83,0] = "c§SMITH1/20§Q§468"
$data[83,1] = 9.51
$data[83,2] = "n"
$data[83,3] = "y"
$data[83,4] = "n"

$data[84,0] = "c§SMITH1/20§Q§468"
$data[84,1] = 6.61
$data[84,2] = "y"
$data[84,3] = "y"
$data[84,4] = "n"

$data[85,0] = "c§JARI1/1/80§Q§173"
$data[85,1] = 12.25
$data[85,2] = "n"
$data[85,4] = "n"

$data[86,0] = "c§DVGM/20§Q§119"
$data[86,1] = 9.19
$data[86,2] = "n"
$data[86,4] = "n"

$data[87,0] = "a§JWBOY/2/80§Q§251"
$data[87,1] = 8.6
$data[87,2] = "n"
$data[87,4] = "n"

$data[88,0] = "b§VRCAN/1/40§Q§487"
$data[88,1] = 8.97
$data[88,2] = "n"
$data[88,3] = "y"
$data[88,4] = "n"

$data[89,0] = "b§VRCAN/1/40§Q§487"
$data[89,1] = 5.98
$data[89,2] = "y"
$data[89,3] = "y"
$data[89,4] = "n"

$data[90,0] = "e§CRSOSO/3/40§Q§54"
$data[90,1] = 8.99
$data[90,2] = "n"
$data[90,4] = "n"

$data[91,0] = "c§SMITH1/20§Q§467"
$data[91,1] = 10.3
$data[91,2] = "n"
$data[91,4] = "n"

$data[92,0] = "b§DVGM/20§Q§113"
$data[92,1] = 10.6
$data[92,2] = "n"
$data[92,4] = "n"

$data[93,0] = "c§MHBUR1/4/20§Q§340"
$data[93,1] = 9.01
$data[93,2] = "n"
$data[93,3] = "y"
$data[93,4] = "n"

$data[94,0] = "c§MHBUR1/4/20§Q§340"
$data[94,1] = 6.28
$data[94,2] = "y"
$data[94,3] = "y"
$data[94,4] = "n"

$data[95,0] = "d§RULEB1/20§N§426"
$data[95,1] = 12.51
$data[95,2] = "n"
$data[95,4] = "n"

$data[96,0] = "c§CWRIC2/X2§Q§95"
$data[96,1] = 10.03
$data[96,2] = "n"
$data[96,4] = "n"

$data[97,0] = "d§JARI1/1/80§Q§179"
$data[97,1] = 9.07
$data[97,2] = "n"
$data[97,4] = "n"

$data[98,0] = "a§DVGM/20§Q§102"
$data[98,1] = 9.6
$data[98,2] = "n"
$data[98,4] = "n"

$data[99,0] = "e§JKFRS1/6/0§V§i122"
$data[99,1] = 14.47
$data[99,2] = "n"
$data[99,4] = "n"

$data[100,0] = "c§NSLJU1/4/0§V§i141"
$data[100,1] = 9.87
$data[100,2] = "n"
$data[100,4] = "n"

$data[101,0] = "e§DVGM/20§Q§131"
$data[101,1] = 17.83
$data[101,2] = "n"
$data[101,4] = "n"

$data[102,0] = "c§JARI1/1/80§Q§177"
$data[102,1] = 9.75
$data[102,2] = "n"
$data[102,4] = "n"

$data[103,0] = "b§SMAKC1/B§N§439"
$data[103,1] = 11.62
$data[103,2] = "n"
$data[103,3] = "y"
$data[103,4] = "n"

$data[104,0] = "b§SMAKC1/B§N§439"
$data[104,1] = 6.24
$data[104,2] = "y"
$data[104,3] = "y"
$data[104,4] = "n"

$data[105,0] = "d§VSGARO1/1/0§Q§543"
$data[105,1] = 10.21
$data[105,2] = "n"
$data[105,4] = "n"

$data[106,0] = "d§KVEDG1/1/80§Q§296"
$data[106,1] = 9.68
$data[106,2] = "n"
$data[106,3] = "y"
$data[106,4] = "n"

$data[107,0] = "d§KVEDG1/1/80§Q§296"
$data[107,1] = 6.45
$data[107,2] = "y"
$data[107,3] = "y"
$data[107,4] = "n"

$data[108,0] = "d§SMAKC1/B§N§442"
$data[108,1] = 9.9
$data[108,2] = "n"
$data[108,4] = "n"

$data[109,0] = "d§VRPET2/3/20§Q§526"
$data[109,1] = 10.1
$data[109,2] = "n"
$data[109,4] = "n"

$data[110,0] = "c§JBCHY1/1/50§Q§235"
$data[110,1] = 10.42
$data[110,2] = "n"
$data[110,4] = "n"

$data[111,0] = "a§SEBRN/1/20§N§i147"
$data[111,1] = 12.56
$data[111,2] = "n"
$data[111,4] = "n"

$data[112,0] = "b§PDVRT1/20§Q§406"
$data[112,1] = 11.36
$data[112,2] = "n"
$data[112,4] = "n"

$data[113,0] = "c§MHBUR1/4/20§Q§349"
$data[113,1] = 12.87
$data[113,2] = "n"
$data[113,4] = "n"

$data[114,0] = "b§KVEDG1/1/80§Q§281"
$data[114,1] = 8.83
$data[114,2] = "n"
$data[114,4] = "n"

$ws.Range("A152:E266").Value = $data
